$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Split-CellIntoTwoRuns($Cell, $FirstText, $SecondText) {
    $cellRange = $Cell.Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $start = $cellRange.Start
    $end = $cellRange.End

    # First run: overwrite the whole original text with the new first chunk.
    $r1 = $d.Range($start, $end)
    $r1.Text = $FirstText

    # Second run: insert the remaining text right after the first run.
    $afterFirst = $start + $FirstText.Length
    $r2 = $d.Range($afterFirst, $afterFirst)
    $r2.InsertAfter($SecondText)
    $r2 = $d.Range($afterFirst, $afterFirst + $SecondText.Length)

    # Toggling Bold off/on forces the engine to keep this as a distinct run
    # instead of silently re-merging it with the adjacent, identically
    # formatted run that was also touched in this editing session.
    $r2.Bold = $false
    $r2.Bold = $true
}

# "Create Trello account and add tasks for the team" row: Duration 2h -> 1h
Split-CellIntoTwoRuns $t.Cell(4, 2) "1" "h"

# "Group Meeting - Gather information..." row: Duration 60 min -> 30 min
Split-CellIntoTwoRuns $t.Cell(5, 2) "30" " min"

# "Setup AWS, Elastic Beanstack, CodePipeline" row: Duration 2h -> 1h
Split-CellIntoTwoRuns $t.Cell(9, 2) "1" "h"
